# #5: cash & deposit done
# Fill in bank/deposit metadata columns (G:M) on the 存款 (deposits) sheet,
# turn row 1 into a proper header row, and fix up row 6's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")

# --- Row 1: convert from a stray data row into the real header row ------
# B1/C1/D1 (bank/deposit_type/currency) already hold the right labels
# (they happen to share text with row 2's first entry), but E1/F1 need to
# become the "owner"/"total" column headers, and G1:M1 are brand new
# metadata headers matching the other sheets' layout.
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"
$ws.Range("G1").Value = "property_category"
$ws.Range("H1").Value = "category"
$ws.Range("I1").Value = "date"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"

# Match the header row's bold/bordered style on the newly added cells.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("G1:M1").PasteSpecial(-4122) | Out-Null

# --- Rows 2-6: fill in the new metadata columns --------------------------
$rows = @(2, 3, 4, 5, 6)
$indexes = @(61, 62, 63, 64, 65)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $ws.Range("G$r").Value = "deposit"
    $ws.Range("H$r").Value = "normal"
    # Force the acquisition date to stay plain text (not auto-converted to
    # a date serial) to match the other property sheets' "date" column.
    $ws.Range("I$r").NumberFormat = "@"
    $ws.Range("I$r").Value = "2012-04-19"
    $ws.Range("J$r").Value = "李慶華"
    $ws.Range("K$r").Value = 607
    $ws.Range("L$r").Value = "tmpe2cb1"
    $ws.Range("M$r").Value = $indexes[$i]
}

# Match rows 2-6's own data style on the new cells (also clears the "@"
# text-format override on column I so it renders like the rest of the row).
$ws.Range("B2").Copy() | Out-Null
$ws.Range("G2:M6").PasteSpecial(-4122) | Out-Null

# --- Row 6 fix-ups: the source index skipped 65, and the amount was ------
# left blank - both get corrected now that the deposit total is known.
$ws.Range("A6").Value = 65
$ws.Range("F6").Value = 543820
